# Update the two e-mail addresses on the "Teste" sheet so that the
# unregistered-email test data points to fresh addresses.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teste")

$ws.Range("D3").Value = "'peterpan30@movie.com"
$ws.Range("D4").Value = "'peterpan31@movie.com"
